$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column P: header "26-jun" plus the day's values for each row.
$ws.Range("P1").Value = "26-jun"

$values = @(
    0,
    15.118505714037614,
    16.042993976342153,
    16.95922399864774,
    0,
    9.7031684073560651,
    7.1668761855984791,
    16.340161153858404,
    12.627995443348777,
    11.147753991943693,
    0,
    12.879422759603091,
    0,
    0,
    14.089944000775855,
    0,
    0
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 16).Value = $values[$i]
}

# Update the selection to match the new active range.
$ws.Range("P2:P18").Select()
